$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.772158333333334
$ws.Range("H2").Value = 14.316475
$ws.Range("I2").Value = 0.2444103987677044
$ws.Range("J2").Value = 0.2444103987677044
$ws.Range("M2").Value = 2.373213
$ws.Range("N2").Value = 7.119638999999999
$ws.Range("O2").Value = 0.6442034269559781
$ws.Range("P2").Value = 0.6442034269559781
$ws.Range("Q2").Value = 11.325348194725
$ws.Range("R2").Value = 101.928133752525
$ws.Range("S2").Value = 0.1574500164698323
$ws.Range("T2").Value = 0.1574500164698323
$ws.Range("G3").Value = 4.772158333333334
$ws.Range("H3").Value = 14.316475
$ws.Range("I3").Value = 0.2444103987677044
$ws.Range("J3").Value = 0.2444103987677044
$ws.Range("O3").Value = 0.2058842823494965
$ws.Range("P3").Value = 0.2058842823494965
$ws.Range("Q3").Value = 3.619526205327778
$ws.Range("R3").Value = 32.57573584795
$ws.Range("S3").Value = 0.05032025954904308
$ws.Range("T3").Value = 0.05032025954904308
$ws.Range("G4").Value = 4.772158333333334
$ws.Range("H4").Value = 14.316475
$ws.Range("I4").Value = 0.2444103987677044
$ws.Range("J4").Value = 0.2444103987677044
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2202126666666667
$ws.Range("N4").Value = 0.6606380000000001
$ws.Range("O4").Value = 0.05977624196638952
$ws.Range("P4").Value = 0.05977624196638952
$ws.Range("Q4").Value = 1.050889712338889
$ws.Range("R4").Value = 9.458007411050001
$ws.Range("S4").Value = 0.01460993513584005
$ws.Range("T4").Value = 0.01460993513584005
$ws.Range("G5").Value = 4.772158333333334
$ws.Range("H5").Value = 14.316475
$ws.Range("I5").Value = 0.2444103987677044
$ws.Range("J5").Value = 0.2444103987677044
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1443573333333333
$ws.Range("N5").Value = 0.433072
$ws.Range("O5").Value = 0.03918547928043534
$ws.Range("P5").Value = 0.03918547928043534
$ws.Range("Q5").Value = 0.6888960512444445
$ws.Range("R5").Value = 6.2000644612
$ws.Range("S5").Value = 0.009577338616834821
$ws.Range("T5").Value = 0.009577338616834821
$ws.Range("G6").Value = 4.772158333333334
$ws.Range("H6").Value = 14.316475
$ws.Range("I6").Value = 0.2444103987677044
$ws.Range("J6").Value = 0.2444103987677044
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09103299999999999
$ws.Range("N6").Value = 0.273099
$ws.Range("O6").Value = 0.0247107067785671
$ws.Range("P6").Value = 0.0247107067785671
$ws.Range("Q6").Value = 0.4344238895583333
$ws.Range("R6").Value = 3.909815006025
$ws.Range("S6").Value = 0.0060395536975814
$ws.Range("T6").Value = 0.006039553697581401
$ws.Range("G7").Value = 4.772158333333334
$ws.Range("H7").Value = 14.316475
$ws.Range("I7").Value = 0.2444103987677044
$ws.Range("J7").Value = 0.2444103987677044
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.09666633333333334
$ws.Range("N7").Value = 0.289999
$ws.Range("O7").Value = 0.02623986266913347
$ws.Range("P7").Value = 0.02623986266913347
$ws.Range("Q7").Value = 0.4613070481694445
$ws.Range("R7").Value = 4.151763433525001
$ws.Range("S7").Value = 0.006413295298572711
$ws.Range("T7").Value = 0.006413295298572711
$ws.Range("G8").Value = 11.54138666666667
$ws.Range("H8").Value = 34.62416
$ws.Range("I8").Value = 0.591102541135077
$ws.Range("J8").Value = 0.591102541135077
$ws.Range("M8").Value = 2.373213
$ws.Range("N8").Value = 7.119638999999999
$ws.Range("O8").Value = 0.6442034269559781
$ws.Range("P8").Value = 0.6442034269559781
$ws.Range("Q8").Value = 27.39016887536
$ws.Range("R8").Value = 246.51151987824
$ws.Range("S8").Value = 0.3807902826816036
$ws.Range("T8").Value = 0.3807902826816036
$ws.Range("G9").Value = 11.54138666666667
$ws.Range("H9").Value = 34.62416
$ws.Range("I9").Value = 0.591102541135077
$ws.Range("J9").Value = 0.591102541135077
$ws.Range("O9").Value = 0.2058842823494965
$ws.Range("P9").Value = 0.2058842823494965
$ws.Range("Q9").Value = 8.753764768035557
$ws.Range("R9").Value = 78.78388291232001
$ws.Range("S9").Value = 0.1216987224765591
$ws.Range("T9").Value = 0.1216987224765591
$ws.Range("G10").Value = 11.54138666666667
$ws.Range("H10").Value = 34.62416
$ws.Range("I10").Value = 0.591102541135077
$ws.Range("J10").Value = 0.591102541135077
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2202126666666667
$ws.Range("N10").Value = 0.6606380000000001
$ws.Range("O10").Value = 0.05977624196638952
$ws.Range("P10").Value = 0.05977624196638952
$ws.Range("Q10").Value = 2.541559534897778
$ws.Range("R10").Value = 22.87403581408
$ws.Range("S10").Value = 0.03533388852583808
$ws.Range("T10").Value = 0.03533388852583808
$ws.Range("G11").Value = 11.54138666666667
$ws.Range("H11").Value = 34.62416
$ws.Range("I11").Value = 0.591102541135077
$ws.Range("J11").Value = 0.591102541135077
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.1443573333333333
$ws.Range("N11").Value = 0.433072
$ws.Range("O11").Value = 0.03918547928043534
$ws.Range("P11").Value = 0.03918547928043534
$ws.Range("Q11").Value = 1.666083802168889
$ws.Range("R11").Value = 14.99475421952
$ws.Range("S11").Value = 0.02316263637826124
$ws.Range("T11").Value = 0.02316263637826124
$ws.Range("G12").Value = 11.54138666666667
$ws.Range("H12").Value = 34.62416
$ws.Range("I12").Value = 0.591102541135077
$ws.Range("J12").Value = 0.591102541135077
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.09103299999999999
$ws.Range("N12").Value = 0.273099
$ws.Range("O12").Value = 0.0247107067785671
$ws.Range("P12").Value = 0.0247107067785671
$ws.Range("Q12").Value = 1.050647052426667
$ws.Range("R12").Value = 9.45582347184
$ws.Range("S12").Value = 0.01460656157005478
$ws.Range("T12").Value = 0.01460656157005478
$ws.Range("G13").Value = 11.54138666666667
$ws.Range("H13").Value = 34.62416
$ws.Range("I13").Value = 0.591102541135077
$ws.Range("J13").Value = 0.591102541135077
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09666633333333334
$ws.Range("N13").Value = 0.289999
$ws.Range("O13").Value = 0.02623986266913347
$ws.Range("P13").Value = 0.02623986266913347
$ws.Range("Q13").Value = 1.115663530648889
$ws.Range("R13").Value = 10.04097177584
$ws.Range("S13").Value = 0.01551044950276024
$ws.Range("T13").Value = 0.01551044950276024
$ws.Range("G14").Value = 3.211640333333333
$ws.Range("H14").Value = 9.634920999999999
$ws.Range("I14").Value = 0.1644870600972187
$ws.Range("J14").Value = 0.1644870600972187
$ws.Range("M14").Value = 2.373213
$ws.Range("N14").Value = 7.119638999999999
$ws.Range("O14").Value = 0.6442034269559781
$ws.Range("P14").Value = 0.6442034269559781
$ws.Range("Q14").Value = 7.621906590390998
$ws.Range("R14").Value = 68.59715931351899
$ws.Range("S14").Value = 0.1059631278045422
$ws.Range("T14").Value = 0.1059631278045422
$ws.Range("G15").Value = 3.211640333333333
$ws.Range("H15").Value = 9.634920999999999
$ws.Range("I15").Value = 0.1644870600972187
$ws.Range("J15").Value = 0.1644870600972187
$ws.Range("O15").Value = 0.2058842823494965
$ws.Range("P15").Value = 0.2058842823494965
$ws.Range("Q15").Value = 2.435924279249111
$ws.Range("R15").Value = 21.923318513242
$ws.Range("S15").Value = 0.03386530032389437
$ws.Range("T15").Value = 0.03386530032389437
$ws.Range("G16").Value = 3.211640333333333
$ws.Range("H16").Value = 9.634920999999999
$ws.Range("I16").Value = 0.1644870600972187
$ws.Range("J16").Value = 0.1644870600972187
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2202126666666667
$ws.Range("N16").Value = 0.6606380000000001
$ws.Range("O16").Value = 0.05977624196638952
$ws.Range("P16").Value = 0.05977624196638952
$ws.Range("Q16").Value = 0.7072438821775555
$ws.Range("R16").Value = 6.365194939597999
$ws.Range("S16").Value = 0.009832418304711399
$ws.Range("T16").Value = 0.009832418304711399
$ws.Range("G17").Value = 3.211640333333333
$ws.Range("H17").Value = 9.634920999999999
$ws.Range("I17").Value = 0.1644870600972187
$ws.Range("J17").Value = 0.1644870600972187
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1443573333333333
$ws.Range("N17").Value = 0.433072
$ws.Range("O17").Value = 0.03918547928043534
$ws.Range("P17").Value = 0.03918547928043534
$ws.Range("Q17").Value = 0.4636238341457777
$ws.Range("R17").Value = 4.172614507312
$ws.Range("S17").Value = 0.006445504285339286
$ws.Range("T17").Value = 0.006445504285339286
$ws.Range("G18").Value = 3.211640333333333
$ws.Range("H18").Value = 9.634920999999999
$ws.Range("I18").Value = 0.1644870600972187
$ws.Range("J18").Value = 0.1644870600972187
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.09103299999999999
$ws.Range("N18").Value = 0.273099
$ws.Range("O18").Value = 0.0247107067785671
$ws.Range("P18").Value = 0.0247107067785671
$ws.Range("Q18").Value = 0.2923652544643333
$ws.Range("R18").Value = 2.631287290178999
$ws.Range("S18").Value = 0.004064591510930915
$ws.Range("T18").Value = 0.004064591510930916
$ws.Range("G19").Value = 3.211640333333333
$ws.Range("H19").Value = 9.634920999999999
$ws.Range("I19").Value = 0.1644870600972187
$ws.Range("J19").Value = 0.1644870600972187
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.09666633333333334
$ws.Range("N19").Value = 0.289999
$ws.Range("O19").Value = 0.02623986266913347
$ws.Range("P19").Value = 0.02623986266913347
$ws.Range("Q19").Value = 0.3104574950087778
$ws.Range("R19").Value = 2.794117455079
$ws.Range("S19").Value = 0.004316117867800522
$ws.Range("T19").Value = 0.004316117867800522
